$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("2024-09-30 21:51:30", "check_availability", "https://example.com", "Checked availability: Selected or default date current date is available for booking.", "2024-09-30", "21:51:30")
    ,@("2024-09-30 21:51:31", "check_availability", "https://example.com", "Failed to check availability: Failed to check availability", "2024-09-30", "21:51:31")
    ,@("2024-09-30 21:51:31", "check_availability", "https://example.com", "Checked availability: No availability for the selected date.", "2024-09-30", "21:51:31")
    ,@("2024-09-30 21:51:31", "check_availability", "https://example.com", "Checked availability: Selected or default date is available for booking.", "2024-09-30", "21:51:31")
    ,@("2024-09-30 21:51:33", "check_availability", "https://example.com", "Failed to check availability: Failed to check availability", "2024-09-30", "21:51:33")
    ,@("2024-09-30 21:54:06", "check_availability", "https://example.com", "Checked availability: Selected or default date current date is available for booking.", "2024-09-30", "21:54:06")
    ,@("2024-09-30 21:54:07", "check_availability", "https://example.com", "Failed to check availability: Failed to check availability", "2024-09-30", "21:54:07")
    ,@("2024-09-30 21:54:07", "check_availability", "https://example.com", "Checked availability: No availability for the selected date.", "2024-09-30", "21:54:07")
    ,@("2024-09-30 21:54:07", "check_availability", "https://example.com", "Checked availability: Selected or default date is available for booking.", "2024-09-30", "21:54:07")
    ,@("2024-09-30 21:54:09", "check_availability", "https://example.com", "Failed to check availability: Failed to check availability", "2024-09-30", "21:54:09")
    ,@("2024-09-30 22:04:35", "check_availability", "https://example.com", "Checked availability: Selected or default date current date is available for booking.", "2024-09-30", "22:04:35")
    ,@("2024-09-30 22:04:36", "check_availability", "https://example.com", "Failed to check availability: Failed to check availability", "2024-09-30", "22:04:36")
    ,@("2024-09-30 22:04:36", "check_availability", "https://example.com", "Checked availability: No availability for the selected date.", "2024-09-30", "22:04:36")
    ,@("2024-09-30 22:07:52", "check_availability", "https://example.com", "Checked availability: Selected or default date current date is available for booking.", "2024-09-30", "22:07:52")
    ,@("2024-09-30 22:07:52", "check_availability", "https://example.com", "Failed to check availability: Failed to check availability", "2024-09-30", "22:07:52")
    ,@("2024-09-30 22:07:52", "check_availability", "https://example.com", "Checked availability: No availability for the selected date.", "2024-09-30", "22:07:52")
    ,@("2024-09-30 22:07:53", "check_availability", "https://example.com", "Checked availability: Selected or default date is available for booking.", "2024-09-30", "22:07:53")
    ,@("2024-09-30 22:07:54", "check_availability", "https://example.com", "Failed to check availability: Failed to check availability", "2024-09-30", "22:07:54")
)

$startRow = 222
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("E$r").ClearFormats()
    $ws.Range("F$r").NumberFormat = "@"
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("F$r").ClearFormats()
}

